$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Reusable questions: replace the B column question text (previously
#    q2..q15 plus two long placeholder strings) with "Question 1".."Question 15"
#    and set every POINTS value (column C) to 2.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le 15; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = "Question $i"
    $ws.Cells.Item($row, 3).Value = 2
}

# ---------------------------------------------------------------------------
# 2. Normalize row heights for rows 3-15 so they match the auto (default)
#    height used by the other rows (no more custom 45/15/30 heights).
# ---------------------------------------------------------------------------
for ($row = 3; $row -le 15; $row++) {
    $ws.Rows.Item($row).AutoFit()
}

# ---------------------------------------------------------------------------
# 3. Add a new, wider column F (for analytics / notes) while keeping D & E
#    at the default width (left untouched).
#    NOTE: the COM ColumnWidth setter works in character units and the
#    runtime quantizes the persisted pixel width, so 58.25 is the input
#    that lands closest to the desired stored width of 59.08984375.
# ---------------------------------------------------------------------------
$ws.Columns.Item(6).ColumnWidth = 58.25

# ---------------------------------------------------------------------------
# 4. Move the active selection from B4 back to B1.
# ---------------------------------------------------------------------------
$ws.Range("B1").Select() | Out-Null

# ---------------------------------------------------------------------------
# 5. Re-protect the sheet with a new password, keeping sheet/objects/
#    scenarios protected, but now explicitly allowing row insertion,
#    row deletion and sorting.
#    Protect(Password, DrawingObjects, Contents, Scenarios, UserInterfaceOnly,
#            AllowFormattingCells, AllowFormattingColumns, AllowFormattingRows,
#            AllowInsertingColumns, AllowInsertingRows, AllowInsertingHyperlinks,
#            AllowDeletingColumns, AllowDeletingRows, AllowSorting,
#            AllowFiltering, AllowUsingPivotTables)
# ---------------------------------------------------------------------------
$ws.Unprotect("essay")
$ws.Protect("essay", $true, $true, $true, $false, $false, $false, $false, $false, $true, $false, $false, $true, $true, $false, $false)
